$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells remain text (avoid Excel auto-converting numeric-looking
# strings like "43.50" or "0.0920" into numbers and dropping trailing zeros).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '42.483.30'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').Value = '2.227.06'
$ws.Range('E3').Value = '  -0.62%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '112.33'
$ws.Range('E5').Value = '  -2.04%  '
$ws.Range('D6').Value = '295.66'
$ws.Range('E6').Value = '  +9.80%  '
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('E8').Value = '  -0.32%  '
$ws.Range('D9').Value = '0.598'
$ws.Range('E9').Value = '  -1.55%  '
$ws.Range('D10').Value = '43.50'
$ws.Range('E10').Value = '  -6.01%  '
$ws.Range('D11').Value = '0.0920'
$ws.Range('E11').Value = '  -0.92%  '
$ws.Range('D12').Value = '54.27'
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('D13').Value = '8.73'
$ws.Range('E13').Value = '  -4.89%  '
$ws.Range('E14').Value = '  +19.96%  '
$ws.Range('E15').Value = '  -1.50%  '
$ws.Range('D16').Value = '14.98'
$ws.Range('E16').Value = '  -2.55%  '
$ws.Range('D17').Value = '2.560.69'
$ws.Range('E17').Value = '  -0.55%  '
$ws.Range('D18').Value = '2.222.59'
$ws.Range('E18').Value = '  -0.94%  '
$ws.Range('D19').Value = '42.430.79'
$ws.Range('E19').Value = '  -1.36%  '
$ws.Range('D20').Value = '7.21'
$ws.Range('E20').Value = '  +6.73%  '
$ws.Range('D22').Value = '73.57'
$ws.Range('E22').Value = '  +2.11%  '
$ws.Range('D23').Value = '3.39'
$ws.Range('E23').Value = '  +15.53%  '
$ws.Range('E24').Value = '  +0.38%  '
$ws.Range('D25').Value = '239.20'
$ws.Range('E25').Value = '  +2.27%  '
$ws.Range('D26').Value = '8.89'
$ws.Range('E26').Value = '  -5.05%  '
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  -1.47%  '
$ws.Range('D28').Value = '11.44'
$ws.Range('E28').Value = '  -6.58%  '
$ws.Range('D29').Value = '2.20'
$ws.Range('E29').Value = '  -1.35%  '
$ws.Range('D30').Value = '175.42'
$ws.Range('E30').Value = '  +0.90%  '
$ws.Range('D31').Value = '37.26'
$ws.Range('E31').Value = '  -8.52%  '
$ws.Range('D32').Value = '21.62'
$ws.Range('E32').Value = '  +2.16%  '
$ws.Range('E33').Value = '  -4.56%  '
$ws.Range('D34').Value = '0.0879'
$ws.Range('E34').Value = '  -3.04%  '
$ws.Range('D35').Value = '5.69'
$ws.Range('E35').Value = '  +1.84%  '
$ws.Range('D36').Value = '4.90'
$ws.Range('E36').Value = '  +4.71%  '
$ws.Range('D37').Value = '0.127'
$ws.Range('E37').Value = '  -0.99%  '
$ws.Range('D38').Value = '4.18'
$ws.Range('E38').Value = '  -2.77%  '
$ws.Range('D39').Value = '0.0373'
$ws.Range('E39').Value = '  -0.25%  '
$ws.Range('E40').Value = '  -2.17%  '
$ws.Range('E41').Value = '  -5.92%  '
$ws.Range('D42').Value = '71.29'
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('D43').Value = '0.228'
$ws.Range('E43').Value = '  -2.23%  '
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('D45').Value = '12.29'
$ws.Range('E45').Value = '  -7.19%  '
$ws.Range('E46').Value = '  -2.42%  '
$ws.Range('E47').Value = '  -4.65%  '
$ws.Range('D48').Value = '1.28'
$ws.Range('E48').Value = '  +2.38%  '
$ws.Range('D49').Value = '8.50'
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('D50').Value = '102.08'
$ws.Range('E50').Value = '  +1.80%  '
$ws.Range('E51').Value = '  +4.61%  '
